$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 (current MSK row) to make room for the new "MOB PRE" row
$ws.Rows.Item(6).Insert()

# Update row 2 (AMM)
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 453.632
$ws.Range("D2").Value = 1225.662
$ws.Range("E2").Value = 481
$ws.Range("F2").Value = 38
$ws.Range("G2").Value = 38
$ws.Range("H2").Value = 258
$ws.Range("I2").Value = 2122.9
$ws.Range("J2").Value = -42.2647322059447

# Update row 3 (IPR)
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 131
$ws.Range("D3").Value = 138
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 288
$ws.Range("J3").Value = -52.08333333333333

# Update row 4 (MIG)
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 464
$ws.Range("D4").Value = 492
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 81
$ws.Range("J4").Value = 507.4074074074074

# Update row 5 (MOB)
$ws.Range("B5").Value = 11
$ws.Range("C5").Value = 767
$ws.Range("D5").Value = 1189
$ws.Range("E5").Value = 341
$ws.Range("F5").Value = 12
$ws.Range("G5").Value = 13
$ws.Range("H5").Value = 58
$ws.Range("I5").Value = 1529
$ws.Range("J5").Value = -22.23675604970569

# Set up new row 6 (MOB PRE) - inherit style from A5 for column A
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A6").Value = "MOB PRE"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 744
$ws.Range("D6").Value = 806
$ws.Range("E6").Value = 49
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 1134
$ws.Range("J6").Value = -28.92416225749559

# Update row 7 (MSK) - previously row 6
$ws.Range("A7").Value = "MSK"
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 275
$ws.Range("D7").Value = 299
$ws.Range("E7").Value = 24
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 321
$ws.Range("J7").Value = -6.853582554517135

# Update row 8 (NOT) - previously row 7
$ws.Range("A8").Value = "NOT"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 82
$ws.Range("D8").Value = 93
$ws.Range("E8").Value = 11
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 2
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 134
$ws.Range("J8").Value = -30.59701492537313

# Update row 9 (TEC) - previously row 8
$ws.Range("A9").Value = "TEC"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 502
$ws.Range("D9").Value = 529
$ws.Range("E9").Value = 22
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 7
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 1307
$ws.Range("J9").Value = -59.5256312165264

# Update row 10 (TST) - previously row 9
$ws.Range("A10").Value = "TST"
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 70
$ws.Range("D10").Value = 95
$ws.Range("E10").Value = 25
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 191
$ws.Range("J10").Value = -50.26178010471204

# Update row 11 (VIP) - previously row 10
$ws.Range("A11").Value = "VIP"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 0

# Update row 12 (WLC) - previously row 11
$ws.Range("A12").Value = "WLC"
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 81
$ws.Range("D12").Value = 87
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 130
$ws.Range("J12").Value = -33.07692307692308
